$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 602 entirely ("問題なし" entry), shifting rows 603:613 up to 602:612
$ws.Rows("602:602").Delete()
